$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the placeholder "   -" values from F2 and F3 (kode_keluarga column)
$ws.Range("F2").ClearContents()
$ws.Range("F3").ClearContents()

# Update the active selection to match the author's final cursor position
$ws.Range("G12").Select()
